# Fixed status in due diligence export
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header D1 was "Prequalification status" -- rename it to the shorter "Status".
$ws.Range("D1").Value = "Status"

# Move the active selection to match where the cursor ended up after the edit
# (row 2 under the renamed "Status" column).
$ws.Range("D2").Select()
